$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.169737577438354
$ws.Range("B1").Value = 2.438621520996094
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.363880395889282
$ws.Range("E1").Value = 1.238795518875122
